# Changes of DEV URL configuration
#
# Replaces the tracking-number text values in column C (rows 2-22) and the
# mirrored values in column D (rows 5-7, 13-17) of Sheet1 with a newer batch
# of tracking numbers, as part of a refreshed shared-strings table for the
# DEV Cheetah-processing test fixture.
#
# The new values are written as genuine TEXT cells (shared-string backed,
# <c t="s">), matching the original cell typing, rather than being
# auto-coerced to numbers the way a bare numeric-looking string normally
# would be when assigned straight to Range.Value. We do this by building the
# literal text in a scratch cell via a formula (="320018799382") - which
# Excel always types as text (t="str") - then Copy / PasteSpecial(values)
# that into the destination cell, and finally clearing the scratch cell. The
# destination cell ends up as a plain shared string with no extra
# NumberFormat/quote-prefix styling, exactly like its neighbours.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Out-of-the-way scratch cell (outside the sheet's real A1:AD24 data range).
$scratch = $ws.Range("ZZ1000")

$updates = @(
  @{Row=2;  Val="320018799382"; HasD=$false},
  @{Row=3;  Val="320018799393"; HasD=$false},
  @{Row=4;  Val="320018799420"; HasD=$false},
  @{Row=5;  Val="320018799441"; HasD=$true},
  @{Row=6;  Val="320018799485"; HasD=$true},
  @{Row=7;  Val="320018799500"; HasD=$true},
  @{Row=8;  Val="320018799533"; HasD=$false},
  @{Row=9;  Val="320018799625"; HasD=$false},
  @{Row=10; Val="320018799658"; HasD=$false},
  @{Row=11; Val="320018799670"; HasD=$false},
  @{Row=12; Val="320018799717"; HasD=$false},
  @{Row=13; Val="320018799739"; HasD=$true},
  @{Row=14; Val="320018799761"; HasD=$true},
  @{Row=15; Val="320018799783"; HasD=$true},
  @{Row=16; Val="320018799810"; HasD=$true},
  @{Row=17; Val="320018799831"; HasD=$true},
  @{Row=18; Val="320018799875"; HasD=$false},
  @{Row=19; Val="320018792701"; HasD=$false},
  @{Row=20; Val="320018792734"; HasD=$false},
  @{Row=21; Val="320018792756"; HasD=$false},
  @{Row=22; Val="320018792789"; HasD=$false}
)

foreach ($item in $updates) {
    $row = $item.Row

    # Build the literal digit string as TEXT in the scratch cell.
    $scratch.Formula = "=""" + $item.Val + """"
    $scratch.Copy()

    $cellC = $ws.Cells.Item($row, 3)
    $cellC.PasteSpecial(-4163)   # xlPasteValues

    if ($item.HasD) {
        $cellD = $ws.Cells.Item($row, 4)
        $cellD.PasteSpecial(-4163)   # xlPasteValues
    }
}

$scratch.ClearContents()
$ws.Range("A1").Select()
